$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. "AppendiX" -> "AppendiCES"  (Contents list entry)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("AppendiX", $true, $false, $false, $false, $false, $true, 1, $false, "AppendiCES", 2)

# ------------------------------------------------------------------
# 2. "8. APPENDIX" -> "8. APPENDICES" (heading)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("8. APPENDIX", $true, $false, $false, $false, $false, $true, 1, $false, "8. APPENDICES", 2)

# ------------------------------------------------------------------
# 3. Remove the word " guide" from "...This guide will also..." and
#    move the _GoBack bookmark to the edit point. Word naturally
#    splits the run that is touched by formatting / the cursor, so we
#    nudge the run boundaries to match (toggle + restore a
#    no-op character format over the middle fragment), then delete
#    " guide" and drop a fresh _GoBack bookmark where it used to be.
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("r future development. This", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitStart = $r.Start
$splitEnd = $r.End

$splitRange = $d.Range($splitStart, $splitEnd)
$splitRange.Font.Bold = 1
$splitRange2 = $d.Range($splitStart, $splitEnd)
$splitRange2.Font.Bold = 0

$r2 = $d.Content
$r2.Find.Execute(" guide will also", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$guideStart = $r2.Start

$delRange = $d.Range($guideStart, $guideStart + 6)
$delRange.Delete()

$bmRange = $d.Range($guideStart, $guideStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ------------------------------------------------------------------
# 4. Footer page number field cached value "5" -> "3"
# ------------------------------------------------------------------
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$fr = $ftr.Range
foreach ($f in $fr.Fields) {
    $res = $f.Result
    if ($res.Text -eq "5") {
        $c = $res.Characters.Item(1)
        $c.Text = "3"
    }
}

Write-Output "done"
